$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add new row of data (row 40)
$ws.Cells.Item(40, 1).Value = 45301
$ws.Cells.Item(40, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(40, 2).Value = 3
$ws.Cells.Item(40, 3).Value = "WebApp Graph/Version editable and deleteable"

# Update the selection to match target
$ws.Range("P50").Select()
